$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1093.4286
$ws.Range("I129").Value = 600
$ws.Range("J129").Value = 1131.3846
$ws.Range("K129").Value = 1800
$ws.Range("L129").Value = 3394.1538
$ws.Range("M129").Value = 3200
$ws.Range("N129").Value = -13394.1538

$ws.Range("H137").Value = 1051.6066
$ws.Range("I137").Value = 905.6591
$ws.Range("J137").Value = 1429.3529
$ws.Range("K137").Value = 2716.9773
$ws.Range("L137").Value = 4288.0587
$ws.Range("M137").Value = -166.9773
$ws.Range("N137").Value = -9388.058700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2694.65
$ws.Range("I32").Value = 2722.7283
$ws.Range("J32").Value = 2371.75
$ws.Range("K32").Value = 2722.7283
$ws.Range("L32").Value = 2371.75
$ws.Range("M32").Value = -2435.7283
$ws.Range("N32").Value = -2945.75

$ws.Range("H61").Value = 4607.7417
$ws.Range("I61").Value = 5650.095
$ws.Range("J61").Value = 2418.8
$ws.Range("K61").Value = 5650.095
$ws.Range("L61").Value = 2418.8
$ws.Range("M61").Value = -5438.095
$ws.Range("N61").Value = -2842.8

$ws.Range("H74").Value = 3261.72
$ws.Range("I74").Value = 3501.318
$ws.Range("J74").Value = 1504.6666
$ws.Range("K74").Value = 3501.318
$ws.Range("L74").Value = 1504.6666
$ws.Range("M74").Value = -2627.318
$ws.Range("N74").Value = -3252.6666

$ws.Range("H77").Value = 3261.72
$ws.Range("I77").Value = 3501.318
$ws.Range("J77").Value = 1504.6666
$ws.Range("K77").Value = 17506.59
$ws.Range("L77").Value = 7523.333000000001
$ws.Range("M77").Value = -13138.59
$ws.Range("N77").Value = -16259.333

$ws.Range("H122").Value = 3022.8462
$ws.Range("I122").Value = 2449.625
$ws.Range("K122").Value = 7348.875
$ws.Range("M122").Value = -4898.875

$ws.Range("H136").Value = 4607.7417
$ws.Range("I136").Value = 5650.095
$ws.Range("J136").Value = 2418.8
$ws.Range("K136").Value = 16950.285
$ws.Range("L136").Value = 7256.400000000001
$ws.Range("M136").Value = -14400.285
$ws.Range("N136").Value = -12356.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1344.0476
$ws.Range("I134").Value = 946.9400000000001
$ws.Range("J134").Value = 2871.3845
$ws.Range("K134").Value = 2840.82
$ws.Range("L134").Value = 8614.1535
$ws.Range("M134").Value = -305.8200000000002
$ws.Range("N134").Value = -13684.1535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1954.9714
$ws.Range("I31").Value = 1189.9362
$ws.Range("J31").Value = 3518.3044
$ws.Range("K31").Value = 1189.9362
$ws.Range("L31").Value = 3518.3044
$ws.Range("M31").Value = -894.9362000000001
$ws.Range("N31").Value = -4108.3044

$ws.Range("H34").Value = 1954.9714
$ws.Range("I34").Value = 1189.9362
$ws.Range("J34").Value = 3518.3044
$ws.Range("K34").Value = 1189.9362
$ws.Range("L34").Value = 3518.3044
$ws.Range("M34").Value = -987.9362000000001
$ws.Range("N34").Value = -3922.3044

$ws.Range("H88").Value = 35000
$ws.Range("J88").Value = 35000
$ws.Range("L88").Value = 35000
$ws.Range("N88").Value = -35812

$ws.Range("H91").Value = 35000
$ws.Range("J91").Value = 35000
$ws.Range("L91").Value = 35000
$ws.Range("N91").Value = -37808

$ws.Range("H95").Value = 13031
$ws.Range("J95").Value = 13031
$ws.Range("L95").Value = 13031
$ws.Range("N95").Value = -18523

$ws.Range("H122").Value = 3158.8572
$ws.Range("I122").Value = 1978
$ws.Range("J122").Value = 4733.3335
$ws.Range("K122").Value = 5934
$ws.Range("L122").Value = 14200.0005
$ws.Range("M122").Value = -3484
$ws.Range("N122").Value = -19100.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1000.4
$ws.Range("I86").Value = 1002
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 3006
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1820
$ws.Range("N86").Value = -5372

$ws.Range("H89").Value = 1000.4
$ws.Range("I89").Value = 1002
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 9018
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -3090
$ws.Range("N89").Value = -20856

$ws.Range("H92").Value = 1073.9286
$ws.Range("I92").Value = 1139
$ws.Range("J92").Value = 1047.9
$ws.Range("K92").Value = 3417
$ws.Range("L92").Value = 3143.7
$ws.Range("M92").Value = -2169
$ws.Range("N92").Value = -5639.700000000001

$ws.Range("H115").Value = 2739.611
$ws.Range("I115").Value = 783.6
$ws.Range("J115").Value = 3491.923
$ws.Range("K115").Value = 2350.8
$ws.Range("L115").Value = 10475.769
$ws.Range("M115").Value = -1175.8
$ws.Range("N115").Value = -12825.769

$ws.Range("H122").Value = 1103.2
$ws.Range("I122").Value = 405.42856
$ws.Range("J122").Value = 2149.8572
$ws.Range("K122").Value = 3648.85704
$ws.Range("L122").Value = 19348.7148
$ws.Range("M122").Value = -1198.85704
$ws.Range("N122").Value = -24248.7148

$ws.Range("H131").Value = 6702.1055
$ws.Range("J131").Value = 10156.667
$ws.Range("L131").Value = 30470.001
$ws.Range("N131").Value = -40550.001

$ws.Range("H134").Value = 1455.6538
$ws.Range("I134").Value = 1071.1333
$ws.Range("J134").Value = 1980
$ws.Range("K134").Value = 3213.3999
$ws.Range("L134").Value = 5940
$ws.Range("M134").Value = 1856.6001
$ws.Range("N134").Value = -16080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 23500
$ws.Range("J95").Value = 23500
$ws.Range("L95").Value = 23500
$ws.Range("N95").Value = -28992

$ws.Range("H102").Value = 1082
$ws.Range("I102").Value = 1098.4
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1098.4
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 523.5999999999999
$ws.Range("N102").Value = -4244

$ws.Range("H132").Value = 1561.3208
$ws.Range("I132").Value = 1221.3182
$ws.Range("J132").Value = 3223.5557
$ws.Range("K132").Value = 3663.9546
$ws.Range("L132").Value = 9670.667099999999
$ws.Range("M132").Value = -1133.9546
$ws.Range("N132").Value = -14730.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 111928
$ws.Range("I46").Value = 167408.67
$ws.Range("J46").Value = 966.6667
$ws.Range("K46").Value = 167408.67
$ws.Range("L46").Value = 966.6667
$ws.Range("M46").Value = -167220.67
$ws.Range("N46").Value = -1342.6667

$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248

$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240

$ws.Range("H132").Value = 4626.023
$ws.Range("I132").Value = 4594.7847
$ws.Range("K132").Value = 13784.3541
$ws.Range("M132").Value = -11254.3541

$ws.Range("H140").Value = 44638.445
$ws.Range("J140").Value = 44638.445
$ws.Range("L140").Value = 44638.445
$ws.Range("N140").Value = -54998.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8551.333000000001
$ws.Range("J41").Value = 8551.333000000001
$ws.Range("L41").Value = 8551.333000000001
$ws.Range("N41").Value = -9331.333000000001

$ws.Range("H45").Value = 8784.5
$ws.Range("I45").Value = 7569
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 7569
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -7078
$ws.Range("N45").Value = -10982
